# Auto-generated edit script applying the crypto price/volume update diff.
# (cryptos list refresh committed by the scheduled GitHub Actions job)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '39.477.82'
$ws.Range('E2').Value = '  +1.66%  '

$ws.Range('D3').Value = '2.162.61'
$ws.Range('E3').Value = '  +3.41%  '

$ws.Range('E4').Value = '  +0.13%  '

$ws.Range('D5').Value = '''229.48'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.40%  '

$ws.Range('D6').Value = '''0.623'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.14%  '

$ws.Range('D7').Value = '''63.08'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.53%  '

$ws.Range('E8').Value = '  +0.11%  '

$ws.Range('D9').Value = '''0.396'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.11%  '

$ws.Range('D10').Value = '''0.0862'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.67%  '

$ws.Range('D11').Value = '''0.103'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.14%  '

$ws.Range('D12').Value = '''16.08'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +7.31%  '

$ws.Range('D13').Value = '2.484.97'
$ws.Range('E13').Value = '  +3.50%  '

$ws.Range('D14').Value = '''22.24'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.24%  '

$ws.Range('D15').Value = '''0.820'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.76%  '

$ws.Range('D16').Value = '''5.58'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.99%  '

$ws.Range('D17').Value = '2.171.75'
$ws.Range('E17').Value = '  +3.71%  '

$ws.Range('D18').Value = '39.468.43'
$ws.Range('E18').Value = '  +1.83%  '

$ws.Range('D19').Value = '''72.37'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.13%  '

$ws.Range('D20').Value = '''6.15'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.76%  '

$ws.Range('D21').Value = '0.0₃0852'
$ws.Range('E21').Value = '  +1.60%  '

$ws.Range('D22').Value = '''228.64'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.62%  '

$ws.Range('E23').Value = '  +0.00%  '

$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D24').Value = '''2.38'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.54%  '

$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').Value = '''2.35'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.46%  '

$ws.Range('D26').Value = '''9.78'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.56%  '

$ws.Range('D27').Value = '''172.31'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.73%  '

$ws.Range('D28').Value = '''0.138'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.00%  '

$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = '''19.66'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.45%  '

$ws.Range('B30').Value = 'ImmutableX'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D30').Value = '''1.42'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.42%  '

$ws.Range('D31').Value = '''2.58'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +7.67%  '

$ws.Range('E32').Value = '  +0.97%  '

$ws.Range('D33').Value = '''4.66'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.51%  '

$ws.Range('D34').Value = '''4.81'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.28%  '

$ws.Range('D35').Value = '''7.01'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +8.24%  '

$ws.Range('D36').Value = '''0.0623'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.88%  '

$ws.Range('D37').Value = '''2.44'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.29%  '

$ws.Range('D38').Value = '''3.55'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.02%  '

$ws.Range('E39').Value = '  +0.06%  '

$ws.Range('B40').Value = 'Aave'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D40').Value = '''103.81'
$ws.Range('D40').Style = 'Normal'

$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').Value = '''0.0231'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.28%  '

$ws.Range('B42').Value = 'InjectiveProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D42').Value = '''18.09'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.09%  '

$ws.Range('D43').Value = '1.528.68'
$ws.Range('E43').Value = '  -0.92%  '

$ws.Range('E44').Value = '  +5.22%  '

$ws.Range('B45').Value = 'Cronos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D45').Value = '''0.0932'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.16%  '

$ws.Range('B46').Value = 'ARBITRUM'
$ws.Range('C46').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D46').Value = '''1.11'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +6.90%  '

$ws.Range('E47').Value = '  -0.10%  '

$ws.Range('B48').Value = 'FTXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D48').Value = '''4.26'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.15%  '

$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').Value = '''7.78'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.61%  '

$ws.Range('D50').Value = '2.369.18'
$ws.Range('E50').Value = '  +3.50%  '

$ws.Range('E51').Value = '  +0.41%  '
